$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 337
$ws.Range("D337").Value = 44615
$ws.Range("K337").Value = 'Abate Fettel'
$ws.Range("L337").Value = 'Primera'
$ws.Range("M337").Value = 200
$ws.Range("N337").Value = 9000
$ws.Range("O337").Value = 10000
$ws.Range("P337").Value = 9500
$ws.Range("S337").Value = 594

# Row 338
$ws.Range("D338").Value = 44615
$ws.Range("K338").Value = 'Abate Fettel'
$ws.Range("L338").Value = 'Segunda'
$ws.Range("M338").Value = 100
$ws.Range("N338").Value = 8000
$ws.Range("O338").Value = 8000
$ws.Range("P338").Value = 8000
$ws.Range("S338").Value = 500

# Row 339
$ws.Range("D339").Value = 44258
$ws.Range("K339").Value = 'Bartlett de verano'
$ws.Range("L339").Value = 'Primera'
$ws.Range("M339").Value = 200
$ws.Range("N339").Value = 9000
$ws.Range("O339").Value = 10000
$ws.Range("P339").Value = 9500
$ws.Range("S339").Value = 594

# Row 340
$ws.Range("D340").Value = 44258
$ws.Range("K340").Value = 'Bartlett de verano'
$ws.Range("L340").Value = 'Segunda'
$ws.Range("M340").Value = 100
$ws.Range("N340").Value = 8000
$ws.Range("O340").Value = 8000
$ws.Range("P340").Value = 8000
$ws.Range("S340").Value = 500

# Row 341
$ws.Range("D341").Value = 44390
$ws.Range("K341").Value = 'Packham''s Triumph'
$ws.Range("L341").Value = 'Primera'
$ws.Range("M341").Value = 200
$ws.Range("N341").Value = 8000
$ws.Range("O341").Value = 9000
$ws.Range("P341").Value = 8500
$ws.Range("S341").Value = 531

# Row 342
$ws.Range("D342").Value = 44390
$ws.Range("K342").Value = 'Packham''s Triumph'
$ws.Range("L342").Value = 'Segunda'
$ws.Range("M342").Value = 100
$ws.Range("N342").Value = 7000
$ws.Range("O342").Value = 7000
$ws.Range("P342").Value = 7000
$ws.Range("S342").Value = 438

# Row 343
$ws.Range("D343").Value = 44498
$ws.Range("K343").Value = 'Packham''s Triumph'
$ws.Range("L343").Value = 'Primera'
$ws.Range("M343").Value = 100
$ws.Range("N343").Value = 9000
$ws.Range("O343").Value = 10000
$ws.Range("P343").Value = 9500
$ws.Range("S343").Value = 594

# Row 344
$ws.Range("D344").Value = 44498
$ws.Range("K344").Value = 'Packham''s Triumph'
$ws.Range("L344").Value = 'Segunda'
$ws.Range("M344").Value = 50
$ws.Range("N344").Value = 8000
$ws.Range("O344").Value = 8000
$ws.Range("P344").Value = 8000
$ws.Range("S344").Value = 500

# Row 345
$ws.Range("D345").Value = 44335
$ws.Range("K345").Value = 'Abate Fettel'
$ws.Range("L345").Value = 'Primera'
$ws.Range("M345").Value = 50
$ws.Range("N345").Value = 9000
$ws.Range("O345").Value = 9000
$ws.Range("P345").Value = 9000
$ws.Range("S345").Value = 562

# Row 346
$ws.Range("D346").Value = 44335
$ws.Range("K346").Value = 'Abate Fettel'
$ws.Range("L346").Value = 'Segunda'
$ws.Range("M346").Value = 50
$ws.Range("N346").Value = 8000
$ws.Range("O346").Value = 8000
$ws.Range("P346").Value = 8000
$ws.Range("S346").Value = 500

# Row 347
$ws.Range("D347").Value = 44552
$ws.Range("K347").Value = 'Packham''s Triumph'
$ws.Range("L347").Value = 'Primera'
$ws.Range("M347").Value = 310
$ws.Range("N347").Value = 9000
$ws.Range("O347").Value = 10000
$ws.Range("P347").Value = 9484
$ws.Range("S347").Value = 593

# Row 348
$ws.Range("A348").Value = 11
$ws.Range("B348").Value = 'Vega Monumental Concepción'
$ws.Range("C348").Value = 'Bíobío'
$ws.Range("D348").Value = 44544
$ws.Range("D348").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E348").Value = 8
$ws.Range("F348").Value = 'Fruta'
$ws.Range("G348").Value = 100104
$ws.Range("H348").Value = 'Frutos de pepita'
$ws.Range("I348").Value = 100104005
$ws.Range("J348").Value = 'Pera'
$ws.Range("K348").Value = 'Packham''s Triumph'
$ws.Range("L348").Value = 'Primera'
$ws.Range("M348").Value = 50
$ws.Range("N348").Value = 11000
$ws.Range("O348").Value = 11000
$ws.Range("P348").Value = 11000
$ws.Range("Q348").Value = '$/caja 16 kilos empedrada'
$ws.Range("R348").Value = 'Región de O''Higgins'
$ws.Range("S348").Value = 688
$ws.Range("T348").Value = 16

# Row 349
$ws.Range("A349").Value = 11
$ws.Range("B349").Value = 'Vega Monumental Concepción'
$ws.Range("C349").Value = 'Bíobío'
$ws.Range("D349").Value = 44544
$ws.Range("D349").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E349").Value = 8
$ws.Range("F349").Value = 'Fruta'
$ws.Range("G349").Value = 100104
$ws.Range("H349").Value = 'Frutos de pepita'
$ws.Range("I349").Value = 100104005
$ws.Range("J349").Value = 'Pera'
$ws.Range("K349").Value = 'Packham''s Triumph'
$ws.Range("L349").Value = 'Segunda'
$ws.Range("M349").Value = 50
$ws.Range("N349").Value = 9000
$ws.Range("O349").Value = 9000
$ws.Range("P349").Value = 9000
$ws.Range("Q349").Value = '$/caja 16 kilos empedrada'
$ws.Range("R349").Value = 'Región de O''Higgins'
$ws.Range("S349").Value = 562
$ws.Range("T349").Value = 16
